$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row number, new Price (D) text (optional), new Volume(1h) (E) text (optional).
# "Num" marks D-values that look numeric so Excel/COM would otherwise silently coerce them
# to a Double (losing trailing zeros like "8.00" -> 8, or reformatting "33.30" -> "33.3").
# For those we pre-format the cell as Text ("@") so the literal string round-trips exactly,
# matching how the source data (plain inline-string cells) is meant to be read.
$updates = @(
    @{ Row = 2; D = "62.409.80"; E = "  -2.01%  " },
    @{ Row = 3; D = "3.166.61"; E = "  -3.84%  " },
    @{ Row = 4; E = "  -0.01%  " },
    @{ Row = 5; D = "585.77"; DNum = $true; E = "  -2.75%  " },
    @{ Row = 6; D = "135.27"; DNum = $true; E = "  -4.50%  " },
    @{ Row = 7; E = "  -0.08%  " },
    @{ Row = 8; D = "3.163.12"; E = "  -3.90%  " },
    @{ Row = 9; D = "0.507"; DNum = $true; E = "  -2.28%  " },
    @{ Row = 10; E = "  -5.42%  " },
    @{ Row = 11; D = "5.26"; DNum = $true; E = "  -3.58%  " },
    @{ Row = 12; D = "0.454"; DNum = $true; E = "  -3.04%  " },
    @{ Row = 13; D = "0.0000234"; DNum = $true; E = "  -5.01%  " },
    @{ Row = 14; D = "33.30"; DNum = $true; E = "  -3.35%  " },
    @{ Row = 15; D = "3.689.71"; E = "  -4.01%  " },
    @{ Row = 16; E = "  -1.93%  " },
    @{ Row = 17; D = "3.164.86"; E = "  -3.81%  " },
    @{ Row = 18; D = "62.387.62"; E = "  -2.16%  " },
    @{ Row = 19; D = "6.53"; DNum = $true; E = "  -4.35%  " },
    @{ Row = 20; D = "451.37"; DNum = $true; E = "  -5.68%  " },
    @{ Row = 21; D = "13.95"; DNum = $true; E = "  -0.83%  " },
    @{ Row = 22; D = "0.700"; DNum = $true; E = "  -3.72%  " },
    @{ Row = 23; D = "7.59"; DNum = $true; E = "  -5.32%  " },
    @{ Row = 24; D = "83.40"; DNum = $true; E = "  -0.85%  " },
    @{ Row = 25; D = "13.26"; DNum = $true; E = "  -1.90%  " },
    @{ Row = 26; E = "  +0.01%  " },
    @{ Row = 27; E = "  -0.08%  " },
    @{ Row = 28; E = "  -3.18%  " },
    @{ Row = 29; D = "6.82"; DNum = $true; E = "  -6.04%  " },
    @{ Row = 30; D = "7.71"; DNum = $true; E = "  -4.74%  " },
    @{ Row = 31; D = "2.01"; DNum = $true; E = "  -7.00%  " },
    @{ Row = 32; D = "27.14"; DNum = $true; E = "  -5.50%  " },
    @{ Row = 33; E = "  -1.41%  " },
    @{ Row = 34; D = "2.38"; DNum = $true; E = "  -5.88%  " },
    @{ Row = 35; E = "  -6.28%  " },
    @{ Row = 36; D = "5.93"; DNum = $true; E = "  -0.52%  " },
    @{ Row = 37; D = "51.21"; DNum = $true; E = "  -3.80%  " },
    @{ Row = 38; D = "0.0₃0697"; E = "  -5.28%  " },
    @{ Row = 39; D = "0.0383"; DNum = $true; E = "  -3.80%  " },
    @{ Row = 40; D = "2.75"; DNum = $true; E = "  +0.85%  " },
    @{ Row = 41; D = "399.87"; DNum = $true; E = "  -6.60%  " },
    @{ Row = 42; D = "8.00"; DNum = $true; E = "  -3.85%  " },
    @{ Row = 43; E = "  -0.99%  " },
    @{ Row = 44; D = "2.798.46"; E = "  -8.61%  " },
    @{ Row = 45; D = "0.249"; DNum = $true; E = "  -5.64%  " },
    @{ Row = 46; E = "  +0.00%  " },
    @{ Row = 47; E = "  -2.60%  " },
    @{ Row = 48; D = "35.54"; DNum = $true; E = "  +1.73%  " },
    @{ Row = 49; D = "124.91"; DNum = $true; E = "  -0.36%  " },
    @{ Row = 50; D = "25.25"; DNum = $true; E = "  -3.36%  " },
    @{ Row = 51; E = "  -3.85%  " }
)

foreach ($item in $updates) {
    if ($item.ContainsKey("D")) {
        $dCell = $ws.Cells.Item($item.Row, 4)
        if ($item.ContainsKey("DNum")) {
            # Force text storage so the exact literal (incl. trailing zeros) is kept.
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($item.Row, 5).Value = $item.E
    }
}
